# Apply updates described by the commit: the player_y record for row 3
# (seas_id 30585 / player_id_y) is changed from De'Aaron Fox's 2022-23
# per-game stats to De'Andre Hunter's 2022-23 per-game stats, and the
# player_id_y on row 2 is corrected (4673 -> 4671).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: fix player_id_y ---
$ws.Range("O2").Value = 4671

# --- Row 3: replace De'Aaron Fox per-game stat line with De'Andre Hunter's ---
$ws.Range("P3").Value  = "De'Andre Hunter"   # player_y
$ws.Range("R3").Value  = "SF"                # pos
$ws.Range("T3").Value  = 4                   # experience
$ws.Range("V3").Value  = "ATL"               # tm_y
$ws.Range("W3").Value  = 67                  # g
$ws.Range("X3").Value  = 67                  # gs
$ws.Range("Y3").Value  = 31.7                # mp_per_game
$ws.Range("Z3").Value  = 5.7                 # fg_per_game
$ws.Range("AA3").Value = 12.3                # fga_per_game
$ws.Range("AB3").Value = 0.461               # fg_percent
$ws.Range("AC3").Value = 1.5                 # x3p_per_game
$ws.Range("AD3").Value = 4.3                 # x3pa_per_game
$ws.Range("AE3").Value = 0.35                # x3p_percent
$ws.Range("AF3").Value = 4.2                 # x2p_per_game
$ws.Range("AG3").Value = 8                   # x2pa_per_game
$ws.Range("AH3").Value = 0.521               # x2p_percent
$ws.Range("AI3").Value = 0.522               # e_fg_percent
$ws.Range("AJ3").Value = 2.6                 # ft_per_game
$ws.Range("AK3").Value = 3.1                 # fta_per_game
$ws.Range("AL3").Value = 0.826               # ft_percent
$ws.Range("AM3").Value = 0.7                 # orb_per_game
$ws.Range("AP3").Value = 1.4                 # ast_per_game
$ws.Range("AQ3").Value = 0.5                 # stl_per_game
$ws.Range("AS3").Value = 1.2                 # tov_per_game
$ws.Range("AT3").Value = 3                   # pf_per_game
$ws.Range("AU3").Value = 15.4                # pts_per_game
